$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

$ws.Range("D2").Value = 15328
$ws.Range("E2").Value = 526
$ws.Range("F2").Value = 526
$ws.Range("G2").Value = 598
$ws.Range("H2").Value = 418
$ws.Range("I2").Value = 435
$ws.Range("J2").Value = -17
$ws.Range("K2").Value = 19559
$ws.Range("L2").Value = 11700
$ws.Range("M2").Value = 7859
$ws.Range("N2").Value = 7748
$ws.Range("O2").Value = 111
$ws.Range("P2").Value = 599
$ws.Range("Q2").Value = 856
$ws.Range("R2").Value = -202
$ws.Range("S2").Value = 142
$ws.Range("T2").Value = 1069
$ws.Range("U2").Value = -213
$ws.Range("V2").Value = 8141
$ws.Range("W2").Value = 3.43
$ws.Range("X2").Value = 2.73
$ws.Range("Y2").Value = 5.93
$ws.Range("Z2").Value = 2.22
$ws.Range("AA2").Value = 148.88
$ws.Range("AB2").Value = 902.38
$ws.Range("AC2").Value = 3634
$ws.Range("AD2").Value = 14.92
$ws.Range("AE2").Value = 65641
$ws.Range("AF2").Value = 0.83
$ws.Range("AH2").Value = 0.74
$ws.Range("AI2").Value = 10.85
$ws.Range("AJ2").Value = 11974656
$ws.Range("D3").Value = 16417
$ws.Range("E3").Value = 411
$ws.Range("F3").Value = 411
$ws.Range("G3").Value = 1302
$ws.Range("H3").Value = 987
$ws.Range("I3").Value = 1003
$ws.Range("J3").Value = -16
$ws.Range("K3").Value = 25037
$ws.Range("L3").Value = 17692
$ws.Range("M3").Value = 7345
$ws.Range("N3").Value = 7174
$ws.Range("O3").Value = 171
$ws.Range("P3").Value = 599
$ws.Range("Q3").Value = 631
$ws.Range("R3").Value = -650
$ws.Range("S3").Value = 215
$ws.Range("T3").Value = 1467
$ws.Range("U3").Value = -836
$ws.Range("V3").Value = 13173
$ws.Range("W3").Value = 2.51
$ws.Range("X3").Value = 6.01
$ws.Range("Y3").Value = 13.44
$ws.Range("Z3").Value = 4.43
$ws.Range("AA3").Value = 240.87
$ws.Range("AB3").Value = 1053.31
$ws.Range("AC3").Value = 8373
$ws.Range("AD3").Value = 5.54
$ws.Range("AE3").Value = 60777
$ws.Range("AF3").Value = 0.76
$ws.Range("AG3").Value = 400
$ws.Range("AH3").Value = 0.86
$ws.Range("AI3").Value = 4.71
$ws.Range("AJ3").Value = 11974656
$ws.Range("D4").Value = 17648
$ws.Range("E4").Value = -153
$ws.Range("F4").Value = -153
$ws.Range("G4").Value = 671
$ws.Range("H4").Value = 376
$ws.Range("I4").Value = 396
$ws.Range("J4").Value = -20
$ws.Range("K4").Value = 25376
$ws.Range("L4").Value = 17673
$ws.Range("M4").Value = 7703
$ws.Range("N4").Value = 7553
$ws.Range("O4").Value = 149
$ws.Range("P4").Value = 599
$ws.Range("Q4").Value = 409
$ws.Range("R4").Value = 367
$ws.Range("S4").Value = -841
$ws.Range("T4").Value = 701
$ws.Range("U4").Value = -292
$ws.Range("V4").Value = 13102
$ws.Range("W4").Value = -0.87
$ws.Range("X4").Value = 2.13
$ws.Range("Y4").Value = 5.37
$ws.Range("Z4").Value = 1.49
$ws.Range("AA4").Value = 229.44
$ws.Range("AB4").Value = 1114.18
$ws.Range("AC4").Value = 3304
$ws.Range("AD4").Value = 8.720000000000001
$ws.Range("AE4").Value = 63994
$ws.Range("AF4").Value = 0.45
$ws.Range("AG4").Value = 400
$ws.Range("AH4").Value = 1.39
$ws.Range("AI4").Value = 11.93
$ws.Range("AJ4").Value = 11974656
$ws.Range("D5").Value = 18126
$ws.Range("E5").Value = 216
$ws.Range("F5").Value = 216
$ws.Range("G5").Value = -468
$ws.Range("H5").Value = -470
$ws.Range("I5").Value = -449
$ws.Range("J5").Value = -21
$ws.Range("K5").Value = 24538
$ws.Range("L5").Value = 15846
$ws.Range("M5").Value = 8692
$ws.Range("N5").Value = 7172
$ws.Range("O5").Value = 1520
$ws.Range("P5").Value = 599
$ws.Range("Q5").Value = 319
$ws.Range("R5").Value = 1534
$ws.Range("S5").Value = -2337
$ws.Range("T5").Value = 731
$ws.Range("U5").Value = -412
$ws.Range("V5").Value = 10561
$ws.Range("W5").Value = 1.19
$ws.Range("X5").Value = -2.59
$ws.Range("Y5").Value = -6.1
$ws.Range("Z5").Value = -1.88
$ws.Range("AA5").Value = 182.3
$ws.Range("AB5").Value = 1038.32
$ws.Range("AC5").Value = -3750
$ws.Range("AD5").Value = -6.93
$ws.Range("AE5").Value = 60765
$ws.Range("AF5").Value = 0.43
$ws.Range("AG5").Value = 400
$ws.Range("AH5").Value = 1.54
$ws.Range("AI5").Value = -10.51
$ws.Range("AJ5").Value = 11974656
$ws.Range("D6").Value = 19508
$ws.Range("E6").Value = 421
$ws.Range("F6").Value = 421
$ws.Range("G6").Value = 673
$ws.Range("H6").Value = 456
$ws.Range("I6").Value = 457
$ws.Range("K6").Value = 26614
$ws.Range("L6").Value = 15943
$ws.Range("M6").Value = 10671
$ws.Range("N6").Value = 9165
$ws.Range("P6").Value = 599
$ws.Range("Q6").Value = 1161
$ws.Range("R6").Value = -570
$ws.Range("S6").Value = -62
$ws.Range("T6").Value = 825
$ws.Range("U6").Value = 336
$ws.Range("V6").Value = 10946
$ws.Range("W6").Value = 2.16
$ws.Range("X6").Value = 2.34
$ws.Range("Y6").Value = 5.6
$ws.Range("Z6").Value = 1.78
$ws.Range("AA6").Value = 149.4
$ws.Range("AB6").Value = 1119.59
$ws.Range("AC6").Value = 3818
$ws.Range("AD6").Value = 14.09
$ws.Range("AE6").Value = 77646
$ws.Range("AF6").Value = 0.6899999999999999
$ws.Range("AG6").Value = 500
$ws.Range("AH6").Value = 0.93
$ws.Range("AI6").Value = 12.91
$ws.Range("AJ6").Value = 11974656
$ws.Range("D7").Value = 20854
$ws.Range("E7").Value = 911
$ws.Range("G7").Value = -27
$ws.Range("H7").Value = -32
$ws.Range("I7").Value = -20
$ws.Range("K7").Value = 34437
$ws.Range("L7").Value = 23874
$ws.Range("M7").Value = 10562
$ws.Range("N7").Value = 9064
$ws.Range("P7").Value = 600
$ws.Range("Q7").Value = 1436
$ws.Range("R7").Value = -911
$ws.Range("S7").Value = -252
$ws.Range("T7").Value = 972
$ws.Range("U7").Value = 158
$ws.Range("W7").Value = 4.37
$ws.Range("X7").Value = -0.15
$ws.Range("Y7").Value = -0.22
$ws.Range("Z7").Value = -0.1
$ws.Range("AA7").Value = 226.05
$ws.Range("AC7").Value = -171
$ws.Range("AD7").Value = -170.28
$ws.Range("AE7").Value = 75693
$ws.Range("AF7").Value = 0.38
$ws.Range("AG7").Value = 432
$ws.Range("AH7").Value = 1.49
$ws.Range("AI7").Value = -253.03
$ws.Range("D8").Value = 22384
$ws.Range("E8").Value = 1023
$ws.Range("G8").Value = 102
$ws.Range("H8").Value = 86
$ws.Range("I8").Value = 92
$ws.Range("K8").Value = 34051
$ws.Range("L8").Value = 23416
$ws.Range("M8").Value = 10635
$ws.Range("N8").Value = 9129
$ws.Range("P8").Value = 600
$ws.Range("Q8").Value = 1049
$ws.Range("R8").Value = -926
$ws.Range("S8").Value = -603
$ws.Range("T8").Value = 1033
$ws.Range("U8").Value = 410
$ws.Range("W8").Value = 4.57
$ws.Range("X8").Value = 0.39
$ws.Range("Y8").Value = 1.01
$ws.Range("Z8").Value = 0.25
$ws.Range("AA8").Value = 220.19
$ws.Range("AC8").Value = 767
$ws.Range("AD8").Value = 43.15
$ws.Range("AE8").Value = 76238
$ws.Range("AF8").Value = 0.43
$ws.Range("AG8").Value = 440
$ws.Range("AH8").Value = 1.33
$ws.Range("AI8").Value = 57.36
$ws.Range("D9").Value = 23803
$ws.Range("E9").Value = 1130
$ws.Range("G9").Value = 253
$ws.Range("H9").Value = 209
$ws.Range("I9").Value = 206
$ws.Range("K9").Value = 34142
$ws.Range("L9").Value = 23347
$ws.Range("M9").Value = 10796
$ws.Range("N9").Value = 9250
$ws.Range("P9").Value = 600
$ws.Range("Q9").Value = 1162
$ws.Range("R9").Value = -600
$ws.Range("S9").Value = -567
$ws.Range("T9").Value = 1064
$ws.Range("U9").Value = 425
$ws.Range("W9").Value = 4.75
$ws.Range("X9").Value = 0.88
$ws.Range("Y9").Value = 2.24
$ws.Range("Z9").Value = 0.61
$ws.Range("AA9").Value = 216.26
$ws.Range("AC9").Value = 1716
$ws.Range("AD9").Value = 19.29
$ws.Range("AE9").Value = 77251
$ws.Range("AF9").Value = 0.43
$ws.Range("AG9").Value = 450
$ws.Range("AH9").Value = 1.36
$ws.Range("AI9").Value = 26.22
